$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("Q2").Value = 2.15
$ws.Range("R2").Value = 1.67

# Row 6 updates
$ws.Range("N6").Value = 8.1
$ws.Range("S6").Value = 1.33
$ws.Range("T6").Value = 3.04
